# Fix typo in the "autoreg_class" column for the boxcox/gaussian row:
# "rm_boxcos" -> "rm_boxcox"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("D7").Value = "rm_boxcox"
